# Auto-generated: apply cryptos price/volume update + coin reordering swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.367.19"
$ws.Range("E2").Value = "  -7.20%  "
$ws.Range("D3").Value = "3.526.48"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'388.95"
$ws.Range("E5").Value = "  -7.13%  "
$ws.Range("D6").Value = "'121.64"
$ws.Range("E6").Value = "  -6.58%  "
$ws.Range("D7").Value = "3.516.07"
$ws.Range("E7").Value = "  -2.63%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("E8").Value = "  -12.05%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.672"
$ws.Range("E10").Value = "  -12.01%  "
$ws.Range("E11").Value = "  -24.22%  "
$ws.Range("D12").Value = "'0.0000319"
$ws.Range("E12").Value = "  -26.67%  "
$ws.Range("D13").Value = "'38.26"
$ws.Range("E13").Value = "  -9.66%  "
$ws.Range("D14").Value = "4.090.79"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "'9.05"
$ws.Range("E15").Value = "  -8.09%  "
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "3.509.09"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("D18").Value = "'12.61"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "'18.55"
$ws.Range("E19").Value = "  -7.43%  "
$ws.Range("D20").Value = "63.449.07"
$ws.Range("E20").Value = "  -6.82%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -10.61%  "
$ws.Range("D22").Value = "'388.27"
$ws.Range("E22").Value = "  -15.68%  "
$ws.Range("D23").Value = "'13.74"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").Value = "'80.31"
$ws.Range("E24").Value = "  -10.87%  "
$ws.Range("D25").Value = "'2.83"
$ws.Range("E25").Value = "  -7.48%  "
$ws.Range("E26").Value = "  +10.17%  "
$ws.Range("D27").Value = "'33.22"
$ws.Range("E27").Value = "  -6.38%  "
$ws.Range("D28").Value = "'2.95"
$ws.Range("E28").Value = "  -10.39%  "
$ws.Range("D29").Value = "'8.58"
$ws.Range("E29").Value = "  -14.92%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'11.64"
$ws.Range("E30").Value = "  -5.03%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.57"
$ws.Range("E31").Value = "  -8.14%  "
$ws.Range("D32").Value = "'0.109"
$ws.Range("E32").Value = "  -8.04%  "
$ws.Range("D33").Value = "'6.69"
$ws.Range("E33").Value = "  -7.07%  "
$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.148"
$ws.Range("E35").Value = "  -5.89%  "
$ws.Range("D36").Value = "'36.27"
$ws.Range("E36").Value = "  -11.00%  "
$ws.Range("D37").Value = "'53.64"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("D38").Value = "'0.0432"
$ws.Range("E38").Value = "  -10.42%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'2.65"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D41").Value = "'0.129"
$ws.Range("E41").Value = "  -12.12%  "
$ws.Range("D42").Value = "'25.97"
$ws.Range("E42").Value = "  +20.47%  "
$ws.Range("D43").Value = "'140.45"
$ws.Range("E43").Value = "  -5.50%  "
$ws.Range("D44").Value = "'3.02"
$ws.Range("E44").Value = "  +15.32%  "
$ws.Range("D45").Value = "'0.0₃0597"
$ws.Range("E45").Value = "  -25.83%  "
$ws.Range("D46").Value = "'1.94"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'2.47"
$ws.Range("E47").Value = "  -9.77%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "'3.04"
$ws.Range("E48").Value = "  -6.50%  "
$ws.Range("D49").Value = "'4.02"
$ws.Range("E49").Value = "  -5.64%  "
$ws.Range("D50").Value = "'2.63"
$ws.Range("E50").Value = "  -11.30%  "
$ws.Range("D51").Value = "'0.271"
$ws.Range("E51").Value = "  -9.96%  "
